$d = $word.ActiveDocument

# Remove the leftover "_GoBack" bookmark (artifact the author's edit removed).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Append a new paragraph after the first one containing two runs:
#   "Учусь использовать " + "github"
$p1 = $d.Paragraphs(1)
$insertionPoint = $d.Range($p1.Range.End, $p1.Range.End)

$newParagraphXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:firstLine="709"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Учусь использовать </w:t></w:r><w:r><w:t>github</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newParagraphXml) | Out-Null
